$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 475
$ws.Range("I4").Value = 475
$ws.Range("K4").Value = 475
$ws.Range("M4").Value = -361

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1371
$ws.Range("J2").Value = 674.5
$ws.Range("L2").Value = 674.5
$ws.Range("N2").Value = -900.5
$ws.Range("H61").Value = 4429.2354
$ws.Range("I61").Value = 3968.2727
$ws.Range("K61").Value = 3968.2727
$ws.Range("M61").Value = -3756.2727
$ws.Range("H63").Value = 7585.7896
$ws.Range("I63").Value = 5974.8
$ws.Range("J63").Value = 9375.777
$ws.Range("K63").Value = 5974.8
$ws.Range("L63").Value = 9375.777
$ws.Range("M63").Value = -5288.8
$ws.Range("N63").Value = -10747.777
$ws.Range("H66").Value = 7585.7896
$ws.Range("I66").Value = 5974.8
$ws.Range("J66").Value = 9375.777
$ws.Range("K66").Value = 29874
$ws.Range("L66").Value = 46878.885
$ws.Range("M66").Value = -26442
$ws.Range("N66").Value = -53742.885
$ws.Range("H116").Value = 1371
$ws.Range("J116").Value = 674.5
$ws.Range("L116").Value = 674.5
$ws.Range("N116").Value = -5262.5
$ws.Range("H122").Value = 1929.0526
$ws.Range("I122").Value = 2041
$ws.Range("J122").Value = 1863.75
$ws.Range("K122").Value = 6123
$ws.Range("L122").Value = 5591.25
$ws.Range("M122").Value = -3673
$ws.Range("N122").Value = -10491.25
$ws.Range("H132").Value = 1746.0938
$ws.Range("I132").Value = 1167.35
$ws.Range("J132").Value = 2710.6667
$ws.Range("K132").Value = 3502.05
$ws.Range("L132").Value = 8132.000100000001
$ws.Range("M132").Value = -972.0499999999997
$ws.Range("N132").Value = -13192.0001
$ws.Range("H136").Value = 4429.2354
$ws.Range("I136").Value = 3968.2727
$ws.Range("K136").Value = 11904.8181
$ws.Range("M136").Value = -9354.8181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1371
$ws.Range("J3").Value = 674.5
$ws.Range("L3").Value = 674.5
$ws.Range("N3").Value = -902.5
$ws.Range("H22").Value = 639.2
$ws.Range("I22").Value = 648.5
$ws.Range("K22").Value = 648.5
$ws.Range("M22").Value = -475.5
$ws.Range("H80").Value = 1439.4667
$ws.Range("I80").Value = 1014.25
$ws.Range("J80").Value = 1594.091
$ws.Range("K80").Value = 1014.25
$ws.Range("L80").Value = 1594.091
$ws.Range("M80").Value = -16.25
$ws.Range("N80").Value = -3590.091
$ws.Range("H83").Value = 1439.4667
$ws.Range("I83").Value = 1014.25
$ws.Range("J83").Value = 1594.091
$ws.Range("K83").Value = 5071.25
$ws.Range("L83").Value = 7970.455
$ws.Range("M83").Value = -79.25
$ws.Range("N83").Value = -17954.455
$ws.Range("H94").Value = 1092.2778
$ws.Range("I94").Value = 1130.7333
$ws.Range("K94").Value = 1130.7333
$ws.Range("M94").Value = -679.7333000000001
$ws.Range("H134").Value = 2734.0334
$ws.Range("I134").Value = 2283.04
$ws.Range("K134").Value = 6849.12
$ws.Range("M134").Value = -4314.12

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3626.7144
$ws.Range("I58").Value = 3697.6
$ws.Range("J58").Value = 3449.5
$ws.Range("K58").Value = 3697.6
$ws.Range("L58").Value = 3449.5
$ws.Range("M58").Value = -3494.6
$ws.Range("N58").Value = -3855.5
$ws.Range("H99").Value = 26967.125
$ws.Range("I99").Value = 7058.2
$ws.Range("K99").Value = 7058.2
$ws.Range("M99").Value = -5560.2
$ws.Range("H121").Value = 75163
$ws.Range("J121").Value = 75163
$ws.Range("L121").Value = 75163
$ws.Range("N121").Value = -77783
$ws.Range("H126").Value = 26967.125
$ws.Range("I126").Value = 7058.2
$ws.Range("K126").Value = 21174.6
$ws.Range("M126").Value = -18704.6
$ws.Range("H134").Value = 3108.4707
$ws.Range("J134").Value = 2468.1428
$ws.Range("L134").Value = 7404.428400000001
$ws.Range("N134").Value = -12474.4284
$ws.Range("H136").Value = 3626.7144
$ws.Range("I136").Value = 3697.6
$ws.Range("J136").Value = 3449.5
$ws.Range("K136").Value = 11092.8
$ws.Range("L136").Value = 10348.5
$ws.Range("M136").Value = -8542.799999999999
$ws.Range("N136").Value = -15448.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 2200
$ws.Range("J51").Value = 2190.5
$ws.Range("L51").Value = 6571.5
$ws.Range("N51").Value = -7491.5
$ws.Range("H60").Value = 228.5
$ws.Range("I60").Value = 198.33333
$ws.Range("J60").Value = 241.42857
$ws.Range("K60").Value = 594.99999
$ws.Range("L60").Value = 724.28571
$ws.Range("M60").Value = -343.99999
$ws.Range("N60").Value = -1226.28571
$ws.Range("H94").Value = 4849.5
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("H122").Value = 681.6667
$ws.Range("J122").Value = 689.75
$ws.Range("L122").Value = 6207.75
$ws.Range("N122").Value = -11107.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7587.9
$ws.Range("I70").Value = 4766.3335
$ws.Range("J70").Value = 8797.143
$ws.Range("K70").Value = 4766.3335
$ws.Range("L70").Value = 8797.143
$ws.Range("M70").Value = -4496.3335
$ws.Range("N70").Value = -9337.143
$ws.Range("H73").Value = 7587.9
$ws.Range("I73").Value = 4766.3335
$ws.Range("J73").Value = 8797.143
$ws.Range("K73").Value = 4766.3335
$ws.Range("L73").Value = 8797.143
$ws.Range("M73").Value = -3830.3335
$ws.Range("N73").Value = -10669.143
$ws.Range("H102").Value = 2730.353
$ws.Range("I102").Value = 2754.0908
$ws.Range("K102").Value = 2754.0908
$ws.Range("M102").Value = -1132.0908
$ws.Range("H132").Value = 3913.3333
$ws.Range("J132").Value = 3774.4443
$ws.Range("L132").Value = 11323.3329
$ws.Range("N132").Value = -16383.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2121.9
$ws.Range("I46").Value = 1879.8334
$ws.Range("J46").Value = 2225.6428
$ws.Range("K46").Value = 1879.8334
$ws.Range("L46").Value = 2225.6428
$ws.Range("M46").Value = -1691.8334
$ws.Range("N46").Value = -2601.6428
$ws.Range("H93").Value = 2995.2173
$ws.Range("J93").Value = 4300.6665
$ws.Range("L93").Value = 4300.6665
$ws.Range("N93").Value = -6796.6665
$ws.Range("H122").Value = 4659.3335
$ws.Range("I122").Value = 3188
$ws.Range("J122").Value = 5710.2856
$ws.Range("K122").Value = 9564
$ws.Range("L122").Value = 17130.8568
$ws.Range("M122").Value = -7114
$ws.Range("N122").Value = -22030.8568
$ws.Range("H132").Value = 4521.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4364.6665
$ws.Range("I132").Value = 4637.8
$ws.Range("K132").Value = 13913.4
$ws.Range("M132").Value = -11383.4
